$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mes" column (C) currently stores the month as a plain number (1-12).
# Replace every value in the data range (C5:C84) with the abbreviated
# Spanish month name so the column holds text instead of numbers.
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 5; $row -le 84; $row++) {
    $cell = $ws.Range("C$row")
    $monthNumber = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNumber]
}

# Match the formatting already used by the note rows right below the table
# (B86/B87) for the "Actualización" note in B85.
$ws.Range("B86").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$excel.CutCopyMode = 0
